# 2023-08-19 update next cosmetic
# The sheet holds a pasted snapshot of a query-style report. This edit
# replaces that snapshot wholesale: the row-number column shifts from B
# to A, every other column shifts left by one (a new col_13 appears in
# N), and the data rows (previously 18) collapse down to 16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old snapshot first so stale cells (old column layout, and the
# now-gone rows 17-18) don't linger.
$ws.Range("A1:M18").ClearContents()

# --- Header row ---
$ws.Range("A1").Value = "row_number"
$ws.Range("B1").Value = "col_1"
$ws.Range("C1").Value = "Type to search:"
$ws.Range("D1").Value = "col_3"
$ws.Range("E1").Value = "col_4"
$ws.Range("F1").Value = "col_5"
$ws.Range("G1").Value = "col_6"
$ws.Range("H1").Value = "col_7"
$ws.Range("I1").Value = "col_8"
$ws.Range("J1").Value = "col_9"
$ws.Range("K1").Value = "col_10"
$ws.Range("L1").Value = "col_11"
$ws.Range("M1").Value = "col_12"
$ws.Range("N1").Value = "col_13"

# --- Row 2 ---
$ws.Range("A2").Value = 3
$ws.Range("C2").Value = "Here is the oldest cosmetic in the backlog by the least-recently added creator, along with any CITB redeems. Change the selector to see each in turn."

# --- Row 3 ---
$ws.Range("A3").Value = 4
$ws.Range("C3").Value = "If there is a red chip in the top right, select the first number in the list to see the current cosmetic  vvvvv"

# --- Row 4 ---
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = " "
$ws.Range("C4").Value = "Shop candidates:"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "Select row to get info:"
$ws.Range("G4").Value = 822

# --- Row 5 ---
$ws.Range("A5").Value = 7
$ws.Range("C5").Value = "Preview"
$ws.Range("D5").Value = "Title/notes"
$ws.Range("F5").Value = "Image links"

# --- Row 6 ---
$ws.Range("A6").Value = 8
$ws.Range("D6").Value = "Lelouch`nCode Geass"
$ws.Range("F6").Value = "Full Head"

# --- Row 7 ---
$ws.Range("A7").Value = 9
$ws.Range("F7").Value = "Hat"
$ws.Range("G7").Value = "https://cdn.discordapp.com/attachments/699111007649398865/1048523933063843931/Dunk_Sweatling_Lelouch_HatHair_V2b_210x210.png"

# --- Row 8 ---
$ws.Range("A8").Value = 10
$ws.Range("F8").Value = "Face"

# --- Row 9 ---
$ws.Range("A9").Value = 11
$ws.Range("F9").Value = "Neck"
$ws.Range("G9").Value = "https://cdn.discordapp.com/attachments/699111007649398865/1048523933399404614/Dunk_Sweatling_Lelouch_NeckClothes_V2b_210x210.png"

# --- Row 10 ---
$ws.Range("A10").Value = 12
$ws.Range("F10").Value = "Body"

# --- Row 11 ---
$ws.Range("A11").Value = 15
$ws.Range("C11").Value = "Creator"
$ws.Range("D11").Value = "Twitch ID"
$ws.Range("E11").Value = "How long have we been waiting?"
$ws.Range("F11").Value = "CITB Redeemer(s)"
$ws.Range("G11").Value = "Redeem Notes"

# --- Row 12 ---
$ws.Range("A12").Value = 16
$ws.Range("C12").Value = "omnipotent_0"
$ws.Range("D12").Value = 42256416
$ws.Range("E12").Value = "259 days"

# --- Row 13 ---
$ws.Range("A13").Value = 18
$ws.Range("C13").Value = "Discord alias(es)"
$ws.Range("D13").Value = "Add aliases and IDs to the UserLookup tab"

# --- Row 14 ---
$ws.Range("A14").Value = 19
$ws.Range("C14").Value = "Omnipotent_0"

# --- Row 15 ---
$ws.Range("A15").Value = 20
$ws.Range("L15").Value = "shop item rows"
$ws.Range("M15").Value = "citb user(s)"
$ws.Range("N15").Value = "citb comment"

# --- Row 16 (was row 17 in the old snapshot; old rows 17-18 are dropped) ---
$ws.Range("A16").Value = 21
$ws.Range("L16").Value = 822
